$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Rename KHUMO brand references (G22 edited first, then B22, to match save order) ---
$ws.Range("G22").Value = "Kumho"
$ws.Range("B22").Value = "Kumo Tires"

# --- Populate Country column (C) ---
# Bulk of rows get "USA" first
$usaRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,19,20,21,23)
foreach ($r in $usaRows) {
    $ws.Range("C$r").Value = "USA"
}

# Then the remaining special countries, in this order: Canada, Germany, Netherlands
$ws.Range("C24").Value = "Canada"
$ws.Range("C17").Value = "Germany"
$ws.Range("C18").Value = "Netherlands"

# --- Fix up formatting for C18:C21 so they match the rest of the Country column (style used by C17/C2) ---
$ws.Range("C17").Copy() | Out-Null
$ws.Range("C18:C21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column C width: no longer auto bestFit, custom width instead ---
$ws.Columns.Item(3).ColumnWidth = 9.5

# --- Row 17 grows taller to fit wrapped "Germany" text ---
$ws.Rows.Item(17).RowHeight = 28

# --- Final selection left on C17 ---
$ws.Range("C17").Select()
